# Notion "CHAM_CONG_SOC_TRANG" (attendance) export sheet gained two more
# tracked days ("Ngay 16"/"Ngay 15" and "Ngay 18"/"Ngay 17"), each
# contributing a 3-column "select" block (id/name/color) that Notion's CSV
# -> xlsx sync inserts into the table, shifting everything to their right.
#
# The two insertion points (before the old "Ngay 13" block, and before the
# old "Ngay 8" block) are reproduced with two real column-insert operations
# so every later column shifts exactly like it would in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two 6-column blocks (shift existing data right) -------
# First block: right before the old "Ngay 13" select-id column (EV).
$ws.Columns("EV:FA").Insert()
# Second block: right before the old "Ngay 8" select-id column, which after
# the first insert now sits at FE.
$ws.Columns("FE:FJ").Insert()

# --- 2. Header row for the 4 newly-added "Ngay N.select.*" triples --------
$ws.Range("EV1").Value = "properties.Ngày 16.select.id"
$ws.Range("EW1").Value = "properties.Ngày 16.select.name"
$ws.Range("EX1").Value = "properties.Ngày 16.select.color"
$ws.Range("EY1").Value = "properties.Ngày 15.select.id"
$ws.Range("EZ1").Value = "properties.Ngày 15.select.name"
$ws.Range("FA1").Value = "properties.Ngày 15.select.color"

$ws.Range("FE1").Value = "properties.Ngày 18.select.id"
$ws.Range("FF1").Value = "properties.Ngày 18.select.name"
$ws.Range("FG1").Value = "properties.Ngày 18.select.color"
$ws.Range("FH1").Value = "properties.Ngày 17.select.id"
$ws.Range("FI1").Value = "properties.Ngày 17.select.name"
$ws.Range("FJ1").Value = "properties.Ngày 17.select.color"

# --- 3. Populate the new columns for the rows whose timesheet now spans --
#        the extra days (each is the usual "full day / pink" select value).
$dataRows = @(5, 9, 11, 17)

foreach ($r in $dataRows) {
    $ws.Range("EV" + $r).Value = "DjwF"
    $ws.Range("EW" + $r).Value = "Đầy đủ"
    $ws.Range("EX" + $r).Value = "pink"
    $ws.Range("EY" + $r).Value = "DjwF"
    $ws.Range("EZ" + $r).Value = "Đầy đủ"
    $ws.Range("FA" + $r).Value = "pink"
    $ws.Range("FE" + $r).Value = "DjwF"
    $ws.Range("FF" + $r).Value = "Đầy đủ"
    $ws.Range("FG" + $r).Value = "pink"
    $ws.Range("FH" + $r).Value = "DjwF"
    $ws.Range("FI" + $r).Value = "Đầy đủ"
    $ws.Range("FJ" + $r).Value = "pink"
}

# --- 4. Bump last_edited_time on the rows touched by the re-sync ---------
foreach ($r in $dataRows) {
    $ws.Range("D" + $r).Value = "2024-07-19T11:41:00.000Z"
}
